$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: update existing headers (B1:AB1) and add new headers (AC1:BH1) ---
$ws.Range("B1").Value = "AAAAACACACAACATCACAAAAATAACCA"
$ws.Range("C1").Value = "AAAAACGCGCGACATCGCGAAAATAACCG"
$ws.Range("D1").Value = "AACACTACCCCCGCGCCTCCTCGCACCCG"
$ws.Range("E1").Value = "AACCAAACGCTCCTACTAACCGCGCACCG"
$ws.Range("F1").Value = "AACCACAACAACCTTAACCCTAAACCCCG"
$ws.Range("G1").Value = "AACCACAACGACCTTAACCCTAAACCCCG"
$ws.Range("H1").Value = "AACCTAACCCCGACAACGCAACTACTCCG"
$ws.Range("I1").Value = "ACCACCACAAAACCCTAAAACTTCTCCCG"
$ws.Range("J1").Value = "ACCACCCCAAATCTATTAATCACCCACCG"
$ws.Range("K1").Value = "ACCGCCACAAAACCCTAAAACTTCTCCCG"
$ws.Range("L1").Value = "ACCGCCACAAAACCCTAAAACTTCTCCNG"
$ws.Range("M1").Value = "ACGATCACTCGATCCACGCGTCCTACCCG"
$ws.Range("N1").Value = "CGAAATCCACTAACGTATAACGAAAACCG"
$ws.Range("O1").Value = "CGGACGTGAAGGGGAGGACGGAGGCGCGT"
$ws.Range("P1").Value = "CGGAGTGTTTTTTTGTAATATTTTTTCGC"
$ws.Range("Q1").Value = "CGGCGTAGGTAGGTTCGTACGAAGTCGTA"
$ws.Range("R1").Value = "CGGCGTAGGTAGGTTCGTACGAAGTTGTA"
$ws.Range("S1").Value = "CGGCGTAGGTAGGTTCGTATGAAGTCGTA"
$ws.Range("T1").Value = "CGGGGAGGTTTATTTGGCGGAAGGAGGGG"
$ws.Range("U1").Value = "CGGGGAGGTTTATTTGGTGGAAGGAGGGG"
$ws.Range("V1").Value = "CGGGGCGGTTTCGTCGAGAAAGGGTGGGA"
$ws.Range("W1").Value = "CGGGGGATAAGGCGTGTTTTAGGGACGTG"
$ws.Range("X1").Value = "CGGGGGATAAGGCGTGTTTTAGGGATGTG"
$ws.Range("Y1").Value = "CGGGGGTTTGGGTCGCGTTTTTTCGTTCG"
$ws.Range("Z1").Value = "CGGGGTTAGGGTTTTTTATGTGCGTAGTA"
$ws.Range("AA1").Value = "CGGGTTTTTAGTTTTTTTGTTATGTGGGA"
$ws.Range("AB1").Value = "CGGTTGTTGGGGTGATCGTAGTTCGTAGC"
$ws.Range("AC1").Value = "CGGTTTAGGGGTAGCGTTACGTTTGGGTT"
$ws.Range("AD1").Value = "CGGTTTTTTTGACGTTATGGTTTTAGGTT"
$ws.Range("AE1").Value = "CGNAGTGTTTTTTTGTAATATTTTTTCGC"
$ws.Range("AF1").Value = "CGNGGTTAGGGTTTTTTATGTGCGTAGTA"
$ws.Range("AG1").Value = "CNGGGCGGTTTCGTCGAGAAAGGGTGGGA"
$ws.Range("AH1").Value = "CTAAACCACCAACACACAAAAAACCACCA"
$ws.Range("AI1").Value = "CTAAACCACCAACACACGAAAAACCACCA"
$ws.Range("AJ1").Value = "CTAAACCACCAACGCGCGAAAAACCGCCG"
$ws.Range("AK1").Value = "CTCCCTAAACGAACACGCGAAACCTCCCA"
$ws.Range("AL1").Value = "CTCCCTAAACGAACACGCGAAACCTCCCG"
$ws.Range("AM1").Value = "CTCCCTAAACGAACACGCGAAACCTCNCA"
$ws.Range("AN1").Value = "GACAACCCTTTAACCGCTAACCTAATCCG"
$ws.Range("AO1").Value = "GACGACCCTTTAACCGCTAACCTAATCCG"
$ws.Range("AP1").Value = "GACGACCCTTTAACCGCTAACCTAATNCG"
$ws.Range("AQ1").Value = "TCTATACCCGCGAATCCACTAAAAACCCA"
$ws.Range("AR1").Value = "TGGAGTGTTTTTTTGTAATATTTTTTTGC"
$ws.Range("AS1").Value = "TGGCGTAGGTAGGTTCGTACGAAGTCGTA"
$ws.Range("AT1").Value = "TGGGAGGGGTTGGGACGGGGCGGGGTTCG"
$ws.Range("AU1").Value = "TGGGAGGGGTTGGGATGGGGTGGGGTTTG"
$ws.Range("AV1").Value = "TGGGGAGGTTTATTTGGCGGAAGGAGGGG"
$ws.Range("AW1").Value = "TGGGGAGGTTTATTTGGTGGAAGGAGGGG"
$ws.Range("AX1").Value = "TGGGGGTTTGGGTCGCGTTTTTTCGTTCG"
$ws.Range("AY1").Value = "TGGGGTTAGGGTTTTTTATGTGTGTAGTA"
$ws.Range("AZ1").Value = "TGGGTTTTCGTGTTGTATTAGTTGTTAGT"
$ws.Range("BA1").Value = "TGGGTTTTTAGTTTTTTCGTTACGTGGGA"
$ws.Range("BB1").Value = "TGGGTTTTTAGTTTTTTTGTTATGTGGGA"
$ws.Range("BC1").Value = "TGGGTTTTTGTGTTGTATTAGTTGTTAGT"
$ws.Range("BD1").Value = "TGNAGTGTTTTTTTGTAATATTTTTTTGC"
$ws.Range("BE1").Value = "TGNGAGGGGTTGGGATGGGGTGGGGTTTG"
$ws.Range("BF1").Value = "TNGAGTGTTTTTTTGTAATATTTTTTTGC"
$ws.Range("BG1").Value = "TNGGAGGGGTTGGGATGGGGTGGGGTTTG"
$ws.Range("BH1").Value = "TNGGGAGGTTTATTTGGTGGAAGGAGGGG"

# Apply the same header style (bold + border, like B1:AB1) to the newly added header cells AC1:BH1
$ws.Range("B1").Copy()
$ws.Range("AC1:BH1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Row 2: clear all data values across B2:BH2 to blank (present but empty) cells ---
$ws.Range("B2:BH2").Value = ""
$ws.Range("B2:BH2").Style = "Normal"

